$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 7.533107333333334
$ws.Range("H2").Value = 22.599322
$ws.Range("I2").Value = 0.4772251808959424
$ws.Range("J2").Value = 0.4772251808959424
$ws.Range("M2").Value = 560.2199806666666
$ws.Range("N2").Value = 1680.659942
$ws.Range("O2").Value = 0.6936344353529325
$ws.Range("P2").Value = 0.6936344353529326
$ws.Range("Q2").Value = 4220.197244639925
$ws.Range("R2").Value = 37981.77520175932
$ws.Range("S2").Value = 0.3310198188869581
$ws.Range("T2").Value = 0.3310198188869581

# Row 3 updates
$ws.Range("G3").Value = 7.533107333333334
$ws.Range("H3").Value = 22.599322
$ws.Range("I3").Value = 0.4772251808959424
$ws.Range("J3").Value = 0.4772251808959424
$ws.Range("O3").Value = 0.06994956469466522
$ws.Range("P3").Value = 0.06994956469466522
$ws.Range("Q3").Value = 425.5857915098827
$ws.Range("R3").Value = 3830.272123588944
$ws.Range("S3").Value = 0.03338169366500403
$ws.Range("T3").Value = 0.03338169366500403

# Row 4 updates
$ws.Range("G4").Value = 7.533107333333334
$ws.Range("H4").Value = 22.599322
$ws.Range("I4").Value = 0.4772251808959424
$ws.Range("J4").Value = 0.4772251808959424
$ws.Range("O4").Value = 0.2364159999524024
$ws.Range("P4").Value = 0.2364159999524024
$ws.Range("Q4").Value = 1438.397664153257
$ws.Range("R4").Value = 12945.57897737931
$ws.Range("S4").Value = 0.1128236683439803
$ws.Range("T4").Value = 0.1128236683439803

# Row 5 updates
$ws.Range("G5").Value = 5.009378000000001
$ws.Range("I5").Value = 0.3173459790819593
$ws.Range("J5").Value = 0.3173459790819593
$ws.Range("M5").Value = 560.2199806666666
$ws.Range("N5").Value = 1680.659942
$ws.Range("O5").Value = 0.6936344353529325
$ws.Range("P5").Value = 0.6936344353529326
$ws.Range("Q5").Value = 2806.353646312025
$ws.Range("R5").Value = 25257.18281680823
$ws.Range("S5").Value = 0.2201220990120384
$ws.Range("T5").Value = 0.2201220990120384

# Row 6 updates
$ws.Range("G6").Value = 5.009378000000001
$ws.Range("I6").Value = 0.3173459790819593
$ws.Range("J6").Value = 0.3173459790819593
$ws.Range("O6").Value = 0.06994956469466522
$ws.Range("P6").Value = 0.06994956469466522
$ws.Range("Q6").Value = 283.006733711152
$ws.Range("S6").Value = 0.02219821309438539
$ws.Range("T6").Value = 0.02219821309438539

# Row 7 updates
$ws.Range("G7").Value = 5.009378000000001
$ws.Range("I7").Value = 0.3173459790819593
$ws.Range("J7").Value = 0.3173459790819593
$ws.Range("O7").Value = 0.2364159999524024
$ws.Range("P7").Value = 0.2364159999524024
$ws.Range("Q7").Value = 956.5080245408309
$ws.Range("R7").Value = 8608.572220867478
$ws.Range("S7").Value = 0.07502566697553557
$ws.Range("T7").Value = 0.07502566697553556

# Row 8 updates
$ws.Range("I8").Value = 0.2054288400220983
$ws.Range("J8").Value = 0.2054288400220983
$ws.Range("M8").Value = 560.2199806666666
$ws.Range("N8").Value = 1680.659942
$ws.Range("O8").Value = 0.6936344353529325
$ws.Range("P8").Value = 0.6936344353529326
$ws.Range("Q8").Value = 1816.648113587014
$ws.Range("R8").Value = 16349.83302228312
$ws.Range("S8").Value = 0.142492517453936
$ws.Range("T8").Value = 0.1424925174539361

# Row 9 updates
$ws.Range("I9").Value = 0.2054288400220983
$ws.Range("J9").Value = 0.2054288400220983
$ws.Range("O9").Value = 0.06994956469466522
$ws.Range("P9").Value = 0.06994956469466522
$ws.Range("S9").Value = 0.0143696579352758
$ws.Range("T9").Value = 0.0143696579352758

# Row 10 updates
$ws.Range("I10").Value = 0.2054288400220983
$ws.Range("J10").Value = 0.2054288400220983
$ws.Range("O10").Value = 0.2364159999524024
$ws.Range("P10").Value = 0.2364159999524024
$ws.Range("S10").Value = 0.04856666463288646
$ws.Range("T10").Value = 0.04856666463288646
